# DEV 4.1 - Officer Registration Management
# Mark the officer's registration status as Approved (was Pending) and
# leave the sheet positioned the way the reviewer left it (cell D3
# selected, zoomed to 191%).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Registration Status (column D) for the one registration row moves from
# "Pending" to "Approved".
$ws.Range("D2").Value = "Approved"

# Leave the view the way it was when the edit was made.
$excel.ActiveWindow.Zoom = 191
$ws.Range("D3").Select()
